$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-26 Monday", 2) | Out-Null
$d.Content.Find.Execute("332÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "679÷7=", 2) | Out-Null
$d.Content.Find.Execute("575÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "668÷3=", 2) | Out-Null
$d.Content.Find.Execute("567÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "812÷8=", 2) | Out-Null
$d.Content.Find.Execute("646÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "215÷4=", 2) | Out-Null
$d.Content.Find.Execute("172÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "636÷5=", 2) | Out-Null
$d.Content.Find.Execute("361÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "758÷7=", 2) | Out-Null
$d.Content.Find.Execute("536÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "684÷4=", 2) | Out-Null
$d.Content.Find.Execute("523÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "529÷8=", 2) | Out-Null
$d.Content.Find.Execute("703÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "299÷7=", 2) | Out-Null
$d.Content.Find.Execute("966÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "933÷8=", 2) | Out-Null
$d.Content.Find.Execute("839÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "824÷3=", 2) | Out-Null
$d.Content.Find.Execute("677÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷4=", 2) | Out-Null
$d.Content.Find.Execute("693÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "975÷7=", 2) | Out-Null
$d.Content.Find.Execute("823÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷2=", 2) | Out-Null
$d.Content.Find.Execute("633÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "465÷5=", 2) | Out-Null
$d.Content.Find.Execute("283÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "499÷4=", 2) | Out-Null
$d.Content.Find.Execute("406÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷9=", 2) | Out-Null
$d.Content.Find.Execute("734÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "892÷3=", 2) | Out-Null
$d.Content.Find.Execute("460÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "538÷4=", 2) | Out-Null
$d.Content.Find.Execute("817÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "619÷9=", 2) | Out-Null
$d.Content.Find.Execute("687÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "499÷3=", 2) | Out-Null
$d.Content.Find.Execute("726÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "353÷5=", 2) | Out-Null
$d.Content.Find.Execute("972÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "455÷6=", 2) | Out-Null
$d.Content.Find.Execute("567÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "502÷7=", 2) | Out-Null
$d.Content.Find.Execute("547÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷6=", 2) | Out-Null
